# Revert "Merge branch 'origin/yetendra' into 'master'" (MR !14)
#
# That merge had split the single "LeaveBalance_48EmployeeCreation" test-data
# row into four rows (employee ranges 1-25 / 26-41 / 83-100 / 100-123),
# pushing every subsequent row down by three. This reverts that: the four
# split rows collapse back into the original single row, the trailing rows
# shift back up, and the TCID numbering in column A is restored to its
# original, contiguous sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three extra rows created by splitting "LeaveBalance_48EmployeeCreation"
# (original row 37 stays; the 38:40 fragments go away).
$ws.Rows("38:40").Delete()

# Remove the three now-duplicated trailing rows (content already shifted up
# into 38:44 by the delete above).
$ws.Rows("45:47").Delete()

# Restore row 37's description to the un-suffixed class name.
$ws.Range("D37").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation"

# Restore contiguous TCID numbering for rows 37-44.
$ws.Range("A37").Value = "42"
$ws.Range("A38").Value = "43"
$ws.Range("A39").Value = "44"
$ws.Range("A40").Value = "45"
$ws.Range("A41").Value = "46"
$ws.Range("A42").Value = "47"
$ws.Range("A43").Value = "48"
$ws.Range("A44").Value = "49"

# Restore the (stale) selection left over from the edit session.
$ws.Range("B56").Select()
